$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"3.015833333333333"
$ws.Cells.Item(2, 8).Value = [double]"9.047499999999999"
$ws.Cells.Item(2, 9).Value = [double]"0.05376901095572644"
$ws.Cells.Item(2, 10).Value = [double]"0.05376901095572643"
$ws.Cells.Item(2, 13).Value = [double]"35.585194"
$ws.Cells.Item(2, 14).Value = [double]"106.755582"
$ws.Cells.Item(2, 15).Value = [double]"0.9972091466993565"
$ws.Cells.Item(2, 16).Value = [double]"0.9972091466993567"
$ws.Cells.Item(2, 17).Value = [double]"107.3190142383333"
$ws.Cells.Item(2, 18).Value = [double]"965.8711281449999"
$ws.Cells.Item(2, 19).Value = [double]"0.05361894953402831"
$ws.Cells.Item(2, 20).Value = [double]"0.05361894953402832"
$ws.Cells.Item(3, 7).Value = [double]"3.015833333333333"
$ws.Cells.Item(3, 8).Value = [double]"9.047499999999999"
$ws.Cells.Item(3, 9).Value = [double]"0.05376901095572644"
$ws.Cells.Item(3, 10).Value = [double]"0.05376901095572643"
$ws.Cells.Item(3, 13).Value = [double]"0.093901"
$ws.Cells.Item(3, 14).Value = [double]"0.281703"
$ws.Cells.Item(3, 15).Value = [double]"0.002631401590341653"
$ws.Cells.Item(3, 16).Value = [double]"0.002631401590341654"
$ws.Cells.Item(3, 17).Value = [double]"0.2831897658333333"
$ws.Cells.Item(3, 18).Value = [double]"2.5487078925"
$ws.Cells.Item(3, 19).Value = [double]"0.0001414878609399963"
$ws.Cells.Item(3, 20).Value = [double]"0.0001414878609399964"
$ws.Cells.Item(4, 7).Value = [double]"3.015833333333333"
$ws.Cells.Item(4, 8).Value = [double]"9.047499999999999"
$ws.Cells.Item(4, 9).Value = [double]"0.05376901095572644"
$ws.Cells.Item(4, 10).Value = [double]"0.05376901095572643"
$ws.Cells.Item(4, 11).Value = [double]"1"
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.00569"
$ws.Cells.Item(4, 14).Value = [double]"0.01707"
$ws.Cells.Item(4, 15).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(4, 16).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(4, 17).Value = [double]"0.01716009166666667"
$ws.Cells.Item(4, 18).Value = [double]"0.154440825"
$ws.Cells.Item(4, 19).Value = [double]"8.573560758123761E-06"
$ws.Cells.Item(4, 20).Value = [double]"8.573560758123761E-06"
$ws.Cells.Item(5, 9).Value = [double]"0.03577730856453667"
$ws.Cells.Item(5, 10).Value = [double]"0.03577730856453666"
$ws.Cells.Item(5, 13).Value = [double]"35.585194"
$ws.Cells.Item(5, 14).Value = [double]"106.755582"
$ws.Cells.Item(5, 15).Value = [double]"0.9972091466993565"
$ws.Cells.Item(5, 16).Value = [double]"0.9972091466993567"
$ws.Cells.Item(5, 17).Value = [double]"71.40889183191935"
$ws.Cells.Item(5, 18).Value = [double]"642.6800264872741"
$ws.Cells.Item(5, 19).Value = [double]"0.03567745934484119"
$ws.Cells.Item(5, 20).Value = [double]"0.03567745934484119"
$ws.Cells.Item(6, 9).Value = [double]"0.03577730856453667"
$ws.Cells.Item(6, 10).Value = [double]"0.03577730856453666"
$ws.Cells.Item(6, 13).Value = [double]"0.093901"
$ws.Cells.Item(6, 14).Value = [double]"0.281703"
$ws.Cells.Item(6, 15).Value = [double]"0.002631401590341653"
$ws.Cells.Item(6, 16).Value = [double]"0.002631401590341654"
$ws.Cells.Item(6, 17).Value = [double]"0.1884313558023334"
$ws.Cells.Item(6, 18).Value = [double]"1.695882202221"
$ws.Cells.Item(6, 19).Value = [double]"9.414446665486584E-05"
$ws.Cells.Item(6, 20).Value = [double]"9.414446665486584E-05"
$ws.Cells.Item(7, 9).Value = [double]"0.03577730856453667"
$ws.Cells.Item(7, 10).Value = [double]"0.03577730856453666"
$ws.Cells.Item(7, 11).Value = [double]"1"
$ws.Cells.Item(7, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(7, 13).Value = [double]"0.00569"
$ws.Cells.Item(7, 14).Value = [double]"0.01707"
$ws.Cells.Item(7, 15).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(7, 16).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(7, 17).Value = [double]"0.01141813627666667"
$ws.Cells.Item(7, 18).Value = [double]"0.10276322649"
$ws.Cells.Item(7, 19).Value = [double]"5.704753040608584E-06"
$ws.Cells.Item(7, 20).Value = [double]"5.704753040608584E-06"
$ws.Cells.Item(8, 7).Value = [double]"1.732509666666666"
$ws.Cells.Item(8, 8).Value = [double]"5.197528999999999"
$ws.Cells.Item(8, 9).Value = [double]"0.03088875310789786"
$ws.Cells.Item(8, 10).Value = [double]"0.03088875310789786"
$ws.Cells.Item(8, 13).Value = [double]"35.585194"
$ws.Cells.Item(8, 14).Value = [double]"106.755582"
$ws.Cells.Item(8, 15).Value = [double]"0.9972091466993565"
$ws.Cells.Item(8, 16).Value = [double]"0.9972091466993567"
$ws.Cells.Item(8, 17).Value = [double]"61.65169259520866"
$ws.Cells.Item(8, 18).Value = [double]"554.865233356878"
$ws.Cells.Item(8, 19).Value = [double]"0.03080254712933391"
$ws.Cells.Item(8, 20).Value = [double]"0.03080254712933392"
$ws.Cells.Item(9, 7).Value = [double]"1.732509666666666"
$ws.Cells.Item(9, 8).Value = [double]"5.197528999999999"
$ws.Cells.Item(9, 9).Value = [double]"0.03088875310789786"
$ws.Cells.Item(9, 10).Value = [double]"0.03088875310789786"
$ws.Cells.Item(9, 13).Value = [double]"0.093901"
$ws.Cells.Item(9, 14).Value = [double]"0.281703"
$ws.Cells.Item(9, 15).Value = [double]"0.002631401590341653"
$ws.Cells.Item(9, 16).Value = [double]"0.002631401590341654"
$ws.Cells.Item(9, 17).Value = [double]"0.1626843902096666"
$ws.Cells.Item(9, 18).Value = [double]"1.464159511887"
$ws.Cells.Item(9, 19).Value = [double]"8.128071405179311E-05"
$ws.Cells.Item(9, 20).Value = [double]"8.128071405179312E-05"
$ws.Cells.Item(10, 7).Value = [double]"1.732509666666666"
$ws.Cells.Item(10, 8).Value = [double]"5.197528999999999"
$ws.Cells.Item(10, 9).Value = [double]"0.03088875310789786"
$ws.Cells.Item(10, 10).Value = [double]"0.03088875310789786"
$ws.Cells.Item(10, 11).Value = [double]"1"
$ws.Cells.Item(10, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(10, 13).Value = [double]"0.00569"
$ws.Cells.Item(10, 14).Value = [double]"0.01707"
$ws.Cells.Item(10, 15).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(10, 16).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(10, 17).Value = [double]"0.009857980003333331"
$ws.Cells.Item(10, 18).Value = [double]"0.08872182002999998"
$ws.Cells.Item(10, 19).Value = [double]"4.925264512142605E-06"
$ws.Cells.Item(10, 20).Value = [double]"4.925264512142606E-06"
$ws.Cells.Item(11, 7).Value = [double]"49.33364366666667"
$ws.Cells.Item(11, 8).Value = [double]"148.000931"
$ws.Cells.Item(11, 9).Value = [double]"0.879564927371839"
$ws.Cells.Item(11, 10).Value = [double]"0.879564927371839"
$ws.Cells.Item(11, 13).Value = [double]"35.585194"
$ws.Cells.Item(11, 14).Value = [double]"106.755582"
$ws.Cells.Item(11, 15).Value = [double]"0.9972091466993565"
$ws.Cells.Item(11, 16).Value = [double]"0.9972091466993567"
$ws.Cells.Item(11, 17).Value = [double]"1755.547280605205"
$ws.Cells.Item(11, 18).Value = [double]"15799.92552544684"
$ws.Cells.Item(11, 19).Value = [double]"0.877110190691153"
$ws.Cells.Item(11, 20).Value = [double]"0.8771101906911531"
$ws.Cells.Item(12, 7).Value = [double]"49.33364366666667"
$ws.Cells.Item(12, 8).Value = [double]"148.000931"
$ws.Cells.Item(12, 9).Value = [double]"0.879564927371839"
$ws.Cells.Item(12, 10).Value = [double]"0.879564927371839"
$ws.Cells.Item(12, 13).Value = [double]"0.093901"
$ws.Cells.Item(12, 14).Value = [double]"0.281703"
$ws.Cells.Item(12, 15).Value = [double]"0.002631401590341653"
$ws.Cells.Item(12, 16).Value = [double]"0.002631401590341654"
$ws.Cells.Item(12, 17).Value = [double]"4.632478473943666"
$ws.Cells.Item(12, 18).Value = [double]"41.692306265493"
$ws.Cells.Item(12, 19).Value = [double]"0.002314488548694998"
$ws.Cells.Item(12, 20).Value = [double]"0.002314488548694998"
$ws.Cells.Item(13, 7).Value = [double]"49.33364366666667"
$ws.Cells.Item(13, 8).Value = [double]"148.000931"
$ws.Cells.Item(13, 9).Value = [double]"0.879564927371839"
$ws.Cells.Item(13, 10).Value = [double]"0.879564927371839"
$ws.Cells.Item(13, 11).Value = [double]"1"
$ws.Cells.Item(13, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 13).Value = [double]"0.00569"
$ws.Cells.Item(13, 14).Value = [double]"0.01707"
$ws.Cells.Item(13, 15).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(13, 16).Value = [double]"0.0001594517103017434"
$ws.Cells.Item(13, 17).Value = [double]"0.2807084324633333"
$ws.Cells.Item(13, 18).Value = [double]"2.52637589217"
$ws.Cells.Item(13, 19).Value = [double]"0.0001402481319908684"
$ws.Cells.Item(13, 20).Value = [double]"0.0001402481319908685"
